# Batch runner is added
#
# MoveInMoveOut_TC_* regression rows on the "SuiteDetails" sheet had their
# Suite column ("I") switched from the "SmokeSuite" test group to the new
# "Regression" test group, and the sheet's scroll/selection position was
# updated to reflect the newly-added bottom rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SuiteDetails")

# Column I (SuiteType) for rows 52-141 (skipping the blank separator rows
# 70, 90, 117, 137) moves from the "SmokeSuite" group to "Regression".
# (Each contiguous block is assigned separately - a single multi-area
# Range.Value assignment only applies to the first area.)
$ws.Range("I52:I69").Value = "Regression"
$ws.Range("I71:I89").Value = "Regression"
$ws.Range("I91:I116").Value = "Regression"
$ws.Range("I118:I136").Value = "Regression"
$ws.Range("I138:I141").Value = "Regression"

# Rows 99-100 previously carried a slightly different cell style (missing
# bottom border) than the rest of the block; restore the common bottom
# border so they match the rest of the column's formatting.
$ws.Range("I99").Borders.Item(9).LineStyle = 1
$ws.Range("I100").Borders.Item(9).LineStyle = 1

# Update the sheet's active selection / scroll position to the new bottom
# of the table.
$ws.Activate()
$ws.Range("I145").Select()
